$wb = $excel.ActiveWorkbook

# --- Idol_School_Dataset sheet: rename Table1 columns via header cells ---
$ws1 = $wb.Worksheets.Item("Idol_School_Dataset")
$ws1.Range("C1").Value2 = "DOB"
$ws1.Range("A1").Value2 = "Name_Chn"
$ws1.Range("B1").Value2 = "Name_Eng"
$ws1.Range("H1").Value2 = "Ability_Rank"
$ws1.Range("I1").Value2 = "Final_Rank"

# --- Produce_48_Dataset sheet: same header renames (plain headers, no table) ---
$ws2 = $wb.Worksheets.Item("Produce_48_Dataset")
$ws2.Range("C1").Value2 = "DOB"
$ws2.Range("A1").Value2 = "Name_Chn"
$ws2.Range("B1").Value2 = "Name_Eng"
$ws2.Range("G1").Value2 = "Final_Rank"

# --- Update the remembered selections on each sheet ---
# Select sheet1's cell first, then sheet2's cell last so Produce_48_Dataset
# ends up as the active/selected tab (matching the saved workbook state).
[void]$ws1.Range("K10").Select()
[void]$ws2.Range("K9").Select()
